# Traitement DGPR/SRNH/SdcAP/BRIL-19 (Enjeu "Autre")
#
# - Numbers the three unnumbered "DDT 31" comments (rows 29, 39, 40) and the
#   four unnumbered "CEREMA/DTerHdF/ASQT/MET" comments (rows 59, 61, 62, 64)
#   in the "Colonne1" column (B). Column A recalculates automatically via its
#   CONCAT formula.
# - Fills in the "Action" column (I) for row 129 (DGPR/SRNH/SdcAP/BRIL-19)
#   with "intgégré" and clears the pending/"to do" yellow highlight.
# - Leaves the active selection on A129, matching where the author ended up.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Synthese PPR")

# --- Numbering the "Colonne1" column for previously-blank rows ---
$ws.Cells.Item(29, 2).Value2 = 1
$ws.Cells.Item(39, 2).Value2 = 2
$ws.Cells.Item(40, 2).Value2 = 3

$ws.Cells.Item(59, 2).Value2 = 1
$ws.Cells.Item(61, 2).Value2 = 2
$ws.Cells.Item(62, 2).Value2 = 3
$ws.Cells.Item(64, 2).Value2 = 4

# --- Action taken for DGPR/SRNH/SdcAP/BRIL-19 ---
$actionCell = $ws.Cells.Item(129, 9)
$actionCell.Value2 = "intgégré"
# Clear the "pending" yellow highlight now that the action is filled in.
$actionCell.Interior.ColorIndex = -4142
$actionCell.Interior.Pattern = -4142

# --- Leave the selection where the author left it ---
$ws.Activate()
[void]$ws.Range("A129").Select()
